$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells whose new value is a plain numeric-looking string must be forced to
# Text format first, otherwise Excel auto-converts them to numbers (losing
# the exact textual representation used by this price sheet).
$textCells = @('D5', 'D6', 'D8', 'D9', 'D10', 'D11', 'D12', 'D13', 'D14', 'D16', 'D18', 'D20', 'D22', 'D23', 'D25', 'D26', 'D27', 'D28', 'D29', 'D32', 'D33', 'D34', 'D35', 'D38', 'D39', 'D40', 'D42', 'D43', 'D44', 'D45', 'D47', 'D48', 'D49', 'D50', 'D51')
foreach ($ref in $textCells) {
    $ws.Range($ref).NumberFormat = "@"
}

$ws.Range('D2').Value = '36.213.17'
$ws.Range('E2').Value = '  -1.57%  '
$ws.Range('D3').Value = '2.037.17'
$ws.Range('E3').Value = '  -1.98%  '
$ws.Range('E4').Value = '  -0.03%  '
$ws.Range('D5').Value = '244.67'
$ws.Range('E5').Value = '  -0.19%  '
$ws.Range('D6').Value = '0.658'
$ws.Range('E7').Value = '  +0.06%  '
$ws.Range('D8').Value = '54.63'
$ws.Range('E8').Value = '  -0.83%  '
$ws.Range('D9').Value = '59.39'
$ws.Range('E9').Value = '  -0.61%  '
$ws.Range('D10').Value = '0.364'
$ws.Range('E10').Value = '  -0.50%  '
$ws.Range('D11').Value = '0.0738'
$ws.Range('E11').Value = '  -2.78%  '
$ws.Range('D12').Value = '0.105'
$ws.Range('E12').Value = '  -3.89%  '
$ws.Range('D13').Value = '0.908'
$ws.Range('E13').Value = '  +2.70%  '
$ws.Range('D14').Value = '14.27'
$ws.Range('E14').Value = '  -4.77%  '
$ws.Range('D15').Value = '2.343.18'
$ws.Range('E15').Value = '  -1.83%  '
$ws.Range('D16').Value = '5.32'
$ws.Range('E16').Value = '  -3.04%  '
$ws.Range('D17').Value = '2.040.86'
$ws.Range('E17').Value = '  -0.96%  '
$ws.Range('D18').Value = '17.40'
$ws.Range('E18').Value = '  +0.54%  '
$ws.Range('D19').Value = '36.154.57'
$ws.Range('E19').Value = '  -1.58%  '
$ws.Range('D20').Value = '71.16'
$ws.Range('E20').Value = '  -2.26%  '
$ws.Range('D21').Value = '0.0₃0850'
$ws.Range('E21').Value = '  -3.22%  '
$ws.Range('D22').Value = '235.66'
$ws.Range('E22').Value = '  -0.72%  '
$ws.Range('D23').Value = '5.17'
$ws.Range('E23').Value = '  -4.78%  '
$ws.Range('E24').Value = '  +0.09%  '
$ws.Range('D25').Value = '2.34'
$ws.Range('E25').Value = '  -2.56%  '
$ws.Range('D26').Value = '2.26'
$ws.Range('E26').Value = '  +4.43%  '
$ws.Range('D27').Value = '9.29'
$ws.Range('E27').Value = '  -5.51%  '
$ws.Range('D28').Value = '163.57'
$ws.Range('E28').Value = '  -2.18%  '
$ws.Range('D29').Value = '19.84'
$ws.Range('E29').Value = '  -3.46%  '
$ws.Range('E30').Value = '  -1.76%  '
$ws.Range('E31').Value = '  -1.49%  '
$ws.Range('D32').Value = '4.94'
$ws.Range('E32').Value = '  -6.74%  '
$ws.Range('D33').Value = '0.0596'
$ws.Range('E33').Value = '  -2.03%  '
$ws.Range('D34').Value = '4.35'
$ws.Range('E34').Value = '  -6.88%  '
$ws.Range('D35').Value = '0.0904'
$ws.Range('E35').Value = '  +8.48%  '
$ws.Range('E36').Value = '  +0.00%  '
$ws.Range('E37').Value = '  -0.79%  '
$ws.Range('D38').Value = '2.20'
$ws.Range('E38').Value = '  -7.24%  '
$ws.Range('D39').Value = '5.04'
$ws.Range('E39').Value = '  +3.65%  '
$ws.Range('D40').Value = '1.20'
$ws.Range('E40').Value = '  -5.86%  '
$ws.Range('E41').Value = '  +1.89%  '
$ws.Range('D42').Value = '0.0214'
$ws.Range('E42').Value = '  -2.86%  '
$ws.Range('D43').Value = '1.09'
$ws.Range('E43').Value = '  -5.40%  '
$ws.Range('B44').Value = 'Cronos'
$ws.Range('C44').Value = 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
$ws.Range('D44').Value = '0.0902'
$ws.Range('E44').Value = '  -5.20%  '
$ws.Range('B45').Value = 'Aave'
$ws.Range('C45').Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
$ws.Range('D45').Value = '92.44'
$ws.Range('E45').Value = '  -4.11%  '
$ws.Range('D46').Value = '1.397.06'
$ws.Range('E46').Value = '  +3.54%  '
$ws.Range('B47').Value = 'FraxShare'
$ws.Range('C47').Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$ws.Range('D47').Value = '7.46'
$ws.Range('E47').Value = '  +3.78%  '
$ws.Range('B48').Value = 'InjectiveProtocol'
$ws.Range('C48').Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
$ws.Range('D48').Value = '15.55'
$ws.Range('E48').Value = '  -2.82%  '
$ws.Range('D49').Value = '2.95'
$ws.Range('E49').Value = '  +2.05%  '
$ws.Range('D50').Value = '2.26'
$ws.Range('E50').Value = '  -6.81%  '
$ws.Range('D51').Value = '45.87'
$ws.Range('E51').Value = '  +1.45%  '
